$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    # Exclude the trailing end-of-cell marker from the range before assigning text,
    # so only the visible content is replaced (formatting of the run is preserved).
    $sub = $d.Range($r.Start, $r.End - 1)
    $sub.Text = $newText
}

# Row 1
Set-CellText $t 1 1 "77×69=5313"
Set-CellText $t 1 2 "32×11=352"
Set-CellText $t 1 3 "51×97=4947"
Set-CellText $t 1 4 "68×86=5848"
Set-CellText $t 1 5 "17×38=646"

# Row 5
Set-CellText $t 5 1 "23×57=1311"
Set-CellText $t 5 2 "47×74=3478"
Set-CellText $t 5 3 "51×47=2397"
Set-CellText $t 5 4 "15×39=585"
Set-CellText $t 5 5 "34×33=1122"

# Row 10
Set-CellText $t 10 1 "43×32=1376"
Set-CellText $t 10 2 "48×14=672"
Set-CellText $t 10 3 "49×47=2303"
Set-CellText $t 10 4 "46×64=2944"
Set-CellText $t 10 5 "48×57=2736"

# Row 15
Set-CellText $t 15 1 "71×14=994"
Set-CellText $t 15 2 "67×12=804"
Set-CellText $t 15 3 "82×15=1230"
Set-CellText $t 15 4 "92×57=5244"
Set-CellText $t 15 5 "89×55=4895"

# Row 20
Set-CellText $t 20 1 "13×54=702"
Set-CellText $t 20 2 "42×11=462"
Set-CellText $t 20 3 "32×89=2848"
Set-CellText $t 20 4 "17×35=595"
Set-CellText $t 20 5 "61×79=4819"
